# Re-shuffle the values of columns A (Id), I (Antal), Q (Ost) and R (Nord)
# among rows 35-47 on the active worksheet, leaving all other cells
# untouched. The mapping below gives, for every row in that block, which
# row originally carried the (A, I, Q, R) values that must end up there
# after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    35 = 38
    36 = 41
    37 = 36
    38 = 39
    39 = 35
    40 = 43
    41 = 47
    42 = 44
    43 = 42
    44 = 40
    45 = 46
    46 = 45
    47 = 37
}

# Column I ("Antal") is stored as text in this sheet (e.g. "1", "2", ...).
# Pre-format it as text so the values we write below are not silently
# re-interpreted by Excel as numbers.
$ws.Range("I35:I47").NumberFormat = "@"

# Snapshot the original values for columns A, I, Q, R before making any
# changes, since several rows will be overwritten during the loop and we
# still need their pre-edit values as sources later in the loop.
$orig = @{}
foreach ($row in $mapping.Values) {
    if (-not $orig.ContainsKey($row)) {
        $orig[$row] = @{
            A = $ws.Cells.Item($row, 1).Value()
            I = [string]$ws.Cells.Item($row, 9).Value()
            Q = $ws.Cells.Item($row, 17).Value()
            R = $ws.Cells.Item($row, 18).Value()
        }
    }
}

foreach ($targetRow in ($mapping.Keys | Sort-Object)) {
    $sourceRow = $mapping[$targetRow]
    $values = $orig[$sourceRow]

    $ws.Cells.Item($targetRow, 1).Value = $values.A
    $ws.Cells.Item($targetRow, 9).Value = $values.I
    $ws.Cells.Item($targetRow, 17).Value = $values.Q
    $ws.Cells.Item($targetRow, 18).Value = $values.R
}
